$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 381
$ws.Range("I33").Value = 319.64285
$ws.Range("K33").Value = 319.64285
$ws.Range("M33").Value = -90.64285000000001
$ws.Range("H98").Value = 317900.62
$ws.Range("I98").Value = 944.25
$ws.Range("J98").Value = 2008334.6
$ws.Range("K98").Value = 944.25
$ws.Range("L98").Value = 2008334.6
$ws.Range("M98").Value = 553.75
$ws.Range("N98").Value = -2011330.6
$ws.Range("H120").Value = 54078.5
$ws.Range("J120").Value = 54078.5
$ws.Range("L120").Value = 54078.5
$ws.Range("N120").Value = -63754.5
$ws.Range("H122").Value = 317900.62
$ws.Range("I122").Value = 944.25
$ws.Range("J122").Value = 2008334.6
$ws.Range("K122").Value = 2832.75
$ws.Range("L122").Value = 6025003.800000001
$ws.Range("M122").Value = -382.75
$ws.Range("N122").Value = -6029903.800000001
$ws.Range("H132").Value = 1285.4615
$ws.Range("I132").Value = 784.36206
$ws.Range("K132").Value = 2353.08618
$ws.Range("M132").Value = 176.9138199999998
$ws.Range("H140").Value = 69190
$ws.Range("J140").Value = 69190
$ws.Range("L140").Value = 69190
$ws.Range("N140").Value = -79550

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4973.701
$ws.Range("J32").Value = 34999
$ws.Range("L32").Value = 34999
$ws.Range("N32").Value = -35573
$ws.Range("H40").Value = 24999.5
$ws.Range("J40").Value = 24999.5
$ws.Range("L40").Value = 24999.5
$ws.Range("N40").Value = -25351.5
$ws.Range("H45").Value = 1811.4286
$ws.Range("I45").Value = 1905.4546
$ws.Range("K45").Value = 1905.4546
$ws.Range("M45").Value = -1528.4546
$ws.Range("H61").Value = 2958.3845
$ws.Range("I61").Value = 2404.5305
$ws.Range("K61").Value = 2404.5305
$ws.Range("M61").Value = -2192.5305
$ws.Range("H136").Value = 2958.3845
$ws.Range("I136").Value = 2404.5305
$ws.Range("K136").Value = 7213.5915
$ws.Range("M136").Value = -4663.5915

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4382.838
$ws.Range("I20").Value = 3850.926
$ws.Range("K20").Value = 3850.926
$ws.Range("M20").Value = -3603.926

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 20353
$ws.Range("H58").Value = 3164.0977
$ws.Range("I58").Value = 1696.6296
$ws.Range("K58").Value = 1696.6296
$ws.Range("M58").Value = -1493.6296
$ws.Range("H60").Value = 32863.363
$ws.Range("J60").Value = 35721.89
$ws.Range("L60").Value = 35721.89
$ws.Range("N60").Value = -36743.89
$ws.Range("H68").Value = 83999.25
$ws.Range("J68").Value = 83999.25
$ws.Range("L68").Value = 83999.25
$ws.Range("N68").Value = -85497.25
$ws.Range("H71").Value = 83999.25
$ws.Range("J71").Value = 83999.25
$ws.Range("L71").Value = 251997.75
$ws.Range("N71").Value = -259485.75
$ws.Range("H99").Value = 2198.577
$ws.Range("I99").Value = 1979.3334
$ws.Range("K99").Value = 1979.3334
$ws.Range("M99").Value = -481.3334
$ws.Range("H126").Value = 2198.577
$ws.Range("I126").Value = 1979.3334
$ws.Range("K126").Value = 5938.0002
$ws.Range("M126").Value = -3468.0002
$ws.Range("H134").Value = 3063.7896
$ws.Range("I134").Value = 2158
$ws.Range("K134").Value = 6474
$ws.Range("M134").Value = -3939
$ws.Range("H136").Value = 3164.0977
$ws.Range("I136").Value = 1696.6296
$ws.Range("K136").Value = 5089.8888
$ws.Range("M136").Value = -2539.8888

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6085.1816
$ws.Range("I3").Value = 991
$ws.Range("K3").Value = 2973
$ws.Range("M3").Value = -2861
$ws.Range("H9").Value = 1000610.75
$ws.Range("I9").Value = 3650060
$ws.Range("J9").Value = 7067.25
$ws.Range("K9").Value = 10950180
$ws.Range("L9").Value = 21201.75
$ws.Range("M9").Value = -10949956
$ws.Range("N9").Value = -21649.75
$ws.Range("H131").Value = 5325956
$ws.Range("I131").Value = 62500600
$ws.Range("J131").Value = 3354416.8
$ws.Range("K131").Value = 187501800
$ws.Range("L131").Value = 10063250.4
$ws.Range("M131").Value = -187496760
$ws.Range("N131").Value = -10073330.4
$ws.Range("H133").Value = 16672817
$ws.Range("I133").Value = 1239.8572
$ws.Range("K133").Value = 3719.5716
$ws.Range("M133").Value = 1340.4284
$ws.Range("H139").Value = 3963.158
$ws.Range("I139").Value = 2098.375
$ws.Range("J139").Value = 5319.364
$ws.Range("K139").Value = 6295.125
$ws.Range("L139").Value = 15958.092
$ws.Range("M139").Value = -1155.125
$ws.Range("N139").Value = -26238.092
$ws.Range("H140").Value = 3099.6155
$ws.Range("I140").Value = 1438.9445
$ws.Range("K140").Value = 4316.833500000001
$ws.Range("M140").Value = 863.1664999999994

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 23857.166
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 23857.166
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 23857.166
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -29847.166
$ws.Range("H122").Value = 4073.7827
$ws.Range("I122").Value = 3941.9167
$ws.Range("J122").Value = 4217.636
$ws.Range("K122").Value = 11825.7501
$ws.Range("L122").Value = 12652.908
$ws.Range("M122").Value = -9375.750100000001
$ws.Range("N122").Value = -17552.908
$ws.Range("H126").Value = 3085.9355
$ws.Range("I126").Value = 1816.9445
$ws.Range("K126").Value = 5450.833500000001
$ws.Range("M126").Value = -2980.833500000001
$ws.Range("H132").Value = 2581.309
$ws.Range("I132").Value = 2246.1133
$ws.Range("J132").Value = 11464
$ws.Range("K132").Value = 6738.3399
$ws.Range("L132").Value = 34392
$ws.Range("M132").Value = -4208.3399
$ws.Range("N132").Value = -39452
$ws.Range("H140").Value = 66528
$ws.Range("J140").Value = 66528
$ws.Range("L140").Value = 66528
$ws.Range("N140").Value = -76888

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 100000
$ws.Range("I45").Value = 100000
$ws.Range("K45").Value = 100000
$ws.Range("M45").Value = -99593
$ws.Range("H132").Value = 4950.125
$ws.Range("I132").Value = 4049.7827
$ws.Range("K132").Value = 12149.3481
$ws.Range("M132").Value = -9619.348100000001
$ws.Range("H136").Value = 4731.892
$ws.Range("I136").Value = 3636.6667
$ws.Range("K136").Value = 10910.0001
$ws.Range("M136").Value = -8360.000100000001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2604.6072
$ws.Range("I126").Value = 2032.381
$ws.Range("J126").Value = 4321.2856
$ws.Range("K126").Value = 6097.143
$ws.Range("L126").Value = 12963.8568
$ws.Range("M126").Value = -3627.143
$ws.Range("N126").Value = -17903.8568
$ws.Range("H132").Value = 3282.4443
$ws.Range("I132").Value = 2442.5
$ws.Range("K132").Value = 7327.5
$ws.Range("M132").Value = -4797.5
$ws.Range("H136").Value = 4994
$ws.Range("I136").Value = 2600.2
$ws.Range("K136").Value = 7800.599999999999
$ws.Range("M136").Value = -5250.599999999999
